# Updates the "cryptos" price/volume table (Sheet1) with refreshed scrape data.
# Mirrors the GitHub Actions commit "Updated cryptos list ... with GitHub Actions":
# each row keeps its Coin/Link unless noted, Price (D) and Volume 1h % (E) are refreshed.
# Price/Volume cells are stored as *text* (not numbers) in the workbook, so plain
# decimal-looking values are written with a leading apostrophe to force text entry,
# then the cell style is reset to "Normal" so no stray numeric/quote-prefix format sticks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.360.08"
$ws.Range("E2").Value = "  -4.86%  "

$ws.Range("D3").Value = "1.563.05"
$ws.Range("E3").Value = "  -5.16%  "

$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("D5").Value = "'1.001"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.18%  "

$ws.Range("D6").Value = "'290.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.29%  "

$ws.Range("D7").Value = "'0.3713"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.40%  "

$ws.Range("D8").Value = "'49.01"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.08%  "

$ws.Range("D9").Value = "'0.3394"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.33%  "

$ws.Range("D10").Value = "'1.165"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.86%  "

$ws.Range("D11").Value = "'0.07640"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.47%  "

$ws.Range("D12").Value = "'1.003"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.38%  "

$ws.Range("D13").Value = "'21.46"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.98%  "

$ws.Range("E14").Value = "  -4.40%  "

$ws.Range("D15").Value = "'6.916"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.88%  "

$ws.Range("D16").Value = "1.570.33"
$ws.Range("E16").Value = "  -5.14%  "

$ws.Range("E17").Value = "  -7.28%  "

$ws.Range("D18").Value = "'89.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.33%  "

$ws.Range("D19").Value = "'0.06723"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.58%  "

$ws.Range("E20").Value = "  +0.19%  "

$ws.Range("D21").Value = "'6.230"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.12%  "

$ws.Range("D22").Value = "'16.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.31%  "

$ws.Range("D23").Value = "'0.5301"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.40%  "

$ws.Range("D24").Value = "'12.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.89%  "

$ws.Range("D25").Value = "22.361.78"
$ws.Range("E25").Value = "  -4.91%  "

$ws.Range("D26").Value = "'2.402"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.17%  "

$ws.Range("D27").Value = "'2.813"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.62%  "

$ws.Range("D28").Value = "'20.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.32%  "

$ws.Range("D29").Value = "'145.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.96%  "

$ws.Range("D30").Value = "'4.984"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.81%  "

$ws.Range("D31").Value = "'125.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.98%  "

$ws.Range("D32").Value = "1.739.88"
$ws.Range("E32").Value = "  -5.31%  "

$ws.Range("D33").Value = "'6.199"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.76%  "

$ws.Range("B34").Value = "WEMIXTOKEN"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'2.008"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.10%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'1.003"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.11%  "

$ws.Range("E36").Value = "  -10.67%  "

$ws.Range("D37").Value = "'0.08477"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.20%  "

$ws.Range("D38").Value = "'0.02528"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.44%  "

$ws.Range("E39").Value = "  -4.28%  "

$ws.Range("D40").Value = "'5.516"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.05%  "

$ws.Range("D41").Value = "'0.06389"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.11%  "

$ws.Range("D42").Value = "'1.289"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.44%  "

$ws.Range("D43").Value = "'11.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -9.13%  "

$ws.Range("D44").Value = "'0.6342"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.76%  "

$ws.Range("D45").Value = "'14.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.16%  "

$ws.Range("E46").Value = "  +0.21%  "

$ws.Range("D47").Value = "'0.5970"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.60%  "

$ws.Range("E48").Value = "  -4.03%  "

$ws.Range("D49").Value = "'2.093"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.13%  "

$ws.Range("D50").Value = "'1.264"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.68%  "

$ws.Range("D51").Value = "'124.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.05%  "
